$wb = $excel.ActiveWorkbook

# Sheet ALC, row 100 (Leve Item ID 19906)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 10312.857
$ws.Range("I100").Value = 10312.857
$ws.Range("K100").Value = 10312.857
$ws.Range("M100").Value = -9771.857

# Sheet ALC, row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 22312.088
$ws.Range("I132").Value = 22312.088
$ws.Range("K132").Value = 66936.264
$ws.Range("M132").Value = -64406.264

# Sheet ALC, row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 36139.5
$ws.Range("I138").Value = 2739.3333
$ws.Range("K138").Value = 8217.999899999999
$ws.Range("M138").Value = -3077.999899999999

# Sheet ARM, row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19826.473
$ws.Range("I32").Value = 21295.53
$ws.Range("J32").Value = 1096
$ws.Range("K32").Value = 21295.53
$ws.Range("L32").Value = 1096
$ws.Range("M32").Value = -21008.53
$ws.Range("N32").Value = -1670

# Sheet ARM, row 63 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2721.8235
$ws.Range("I63").Value = 2584.8
$ws.Range("K63").Value = 2584.8
$ws.Range("M63").Value = -1898.8

# Sheet ARM, row 66 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2721.8235
$ws.Range("I66").Value = 2584.8
$ws.Range("K66").Value = 12924
$ws.Range("M66").Value = -9492

# Sheet ARM, row 88 (Leve Item ID 12530)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 6797.1
$ws.Range("J88").Value = 8799.429
$ws.Range("L88").Value = 8799.429
$ws.Range("N88").Value = -9611.429

# Sheet ARM, row 91 (Leve Item ID 12530)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 6797.1
$ws.Range("J91").Value = 8799.429
$ws.Range("L91").Value = 8799.429
$ws.Range("N91").Value = -11607.429

# Sheet ARM, row 97 (Leve Item ID 19941)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 995.26086
$ws.Range("I97").Value = 927.619
$ws.Range("J97").Value = 1705.5
$ws.Range("K97").Value = 927.619
$ws.Range("L97").Value = 1705.5
$ws.Range("M97").Value = -431.619
$ws.Range("N97").Value = -2697.5

# Sheet ARM, row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1654.909
$ws.Range("I132").Value = 1245.4445
$ws.Range("K132").Value = 3736.3335
$ws.Range("M132").Value = -1206.3335

# Sheet ARM, row 140 (Leve Item ID 42496)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# Sheet BSM, row 20 (Leve Item ID 14149)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14476.48
$ws.Range("I20").Value = 25835
$ws.Range("K20").Value = 25835
$ws.Range("M20").Value = -25588

# Sheet BSM, row 86 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1570.6471
$ws.Range("I86").Value = 1578.2222
$ws.Range("J86").Value = 1562.125
$ws.Range("K86").Value = 1578.2222
$ws.Range("L86").Value = 1562.125
$ws.Range("M86").Value = -455.2221999999999
$ws.Range("N86").Value = -3808.125

# Sheet BSM, row 89 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1570.6471
$ws.Range("I89").Value = 1578.2222
$ws.Range("J89").Value = 1562.125
$ws.Range("K89").Value = 7891.111
$ws.Range("L89").Value = 7810.625
$ws.Range("M89").Value = -2275.111
$ws.Range("N89").Value = -19042.625

# Sheet BSM, row 94 (Leve Item ID 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2149
$ws.Range("I94").Value = 1688.7778
$ws.Range("J94").Value = 3529.6667
$ws.Range("K94").Value = 1688.7778
$ws.Range("L94").Value = 3529.6667
$ws.Range("M94").Value = -1237.7778
$ws.Range("N94").Value = -4431.6667

# Sheet BSM, row 107 (Leve Item ID 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3028.8857
$ws.Range("I107").Value = 2949.92
$ws.Range("K107").Value = 2949.92
$ws.Range("M107").Value = -1029.92

# Sheet CRP, row 4 (Leve Item ID 3742)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# Sheet CRP, row 74 (Leve Item ID 10636)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 55000
$ws.Range("J74").Value = 55000
$ws.Range("L74").Value = 55000
$ws.Range("N74").Value = -56748

# Sheet CRP, row 77 (Leve Item ID 10636)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 55000
$ws.Range("J77").Value = 55000
$ws.Range("L77").Value = 165000
$ws.Range("N77").Value = -173736

# Sheet CUL, row 36 (Leve Item ID 4732)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

# Sheet CUL, row 92 (Leve Item ID 19841)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 689.25
$ws.Range("I92").Value = 766.8333
$ws.Range("J92").Value = 456.5
$ws.Range("K92").Value = 2300.4999
$ws.Range("L92").Value = 1369.5
$ws.Range("M92").Value = -1052.4999
$ws.Range("N92").Value = -3865.5

# Sheet CUL, row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2823.2188
$ws.Range("I131").Value = 6285
$ws.Range("J131").Value = 2024.3462
$ws.Range("K131").Value = 18855
$ws.Range("L131").Value = 6073.0386
$ws.Range("M131").Value = -13815
$ws.Range("N131").Value = -16153.0386

# Sheet GSM, row 95 (Leve Item ID 18235)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 74492
$ws.Range("J95").Value = 74492
$ws.Range("L95").Value = 74492
$ws.Range("N95").Value = -79984

# Sheet GSM, row 107 (Leve Item ID 27802)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 388.69232
$ws.Range("I107").Value = 104.625
$ws.Range("K107").Value = 104.625
$ws.Range("M107").Value = 1815.375

# Sheet GSM, row 113 (Leve Item ID 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2414.5293
$ws.Range("I113").Value = 1502.5
$ws.Range("K113").Value = 1502.5
$ws.Range("M113").Value = 667.5

# Sheet GSM, row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2816.0435
$ws.Range("I122").Value = 2489
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 7467
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -5017
$ws.Range("N122").Value = -23650

# Sheet LTW, row 2 (Leve Item ID 2631)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2544997.2
$ws.Range("J2").Value = 2544997.2
$ws.Range("L2").Value = 2544997.2
$ws.Range("N2").Value = -2545221.2

# Sheet LTW, row 7 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5916.3335
$ws.Range("I7").Value = 4749.75
$ws.Range("J7").Value = 8249.5
$ws.Range("K7").Value = 4749.75
$ws.Range("L7").Value = 8249.5
$ws.Range("M7").Value = -4637.75
$ws.Range("N7").Value = -8473.5

# Sheet LTW, row 126 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5916.3335
$ws.Range("I126").Value = 4749.75
$ws.Range("J126").Value = 8249.5
$ws.Range("K126").Value = 14249.25
$ws.Range("L126").Value = 24748.5
$ws.Range("M126").Value = -11779.25
$ws.Range("N126").Value = -29688.5

# Sheet LTW, row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2818.44
$ws.Range("I132").Value = 2534.9
$ws.Range("K132").Value = 7604.700000000001
$ws.Range("M132").Value = -5074.700000000001

# Sheet WVR, row 80 (Leve Item ID 10911)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 31000
$ws.Range("I80").Value = 31000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 31000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -30002
$ws.Range("N80").ClearContents()

# Sheet WVR, row 83 (Leve Item ID 10911)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 31000
$ws.Range("I83").Value = 31000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 93000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -88008
$ws.Range("N83").ClearContents()

# Sheet WVR, row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3782.9
$ws.Range("I126").Value = 3214.5293
$ws.Range("K126").Value = 9643.5879
$ws.Range("M126").Value = -7173.5879
